# Modificata gestione DB e aggiunta dei mezzi sulla mappa
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the numeric lat/lon values in row 5 with their textual
# (shared-string) equivalents, matching the new coordinates used for the
# "Piazza San Vigilio" station.
$ws.Range("E5").Value = "45.708509"
$ws.Range("F5").Value = "9.650654"

# Move the active selection from E3 to F5, as recorded in the saved view.
$ws.Range("F5").Select()
